$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# --- Status text update: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"

# Overview mirrors the Status value for each locale in columns E (zh-cn) and F (de-de)
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime updates ---
$wsZh.Range("K2").Value = "2016-09-07 07:09:50"
$wsDe.Range("K2").Value = "2016-09-07 07:09:59"

# --- Error Detail cleared (handback is now in sync, no error) ---
$wsZh.Range("P2").Value = ""
$wsDe.Range("P2").Value = ""

# --- Column width changes ---
# (target OOXML widths are 29.9777047293527 and 13.7470528738839; the COM
# layer snaps ColumnWidth to the nearest 1/6-character pixel grid like real
# Excel, so we pick the input that lands on the closest achievable width.)
$wsOverview.Range("E1").ColumnWidth = 29.166666666666668
$wsOverview.Range("F1").ColumnWidth = 29.166666666666668

$wsZh.Range("C1").ColumnWidth = 29.166666666666668
$wsZh.Range("P1").ColumnWidth = 12.833333333333334

$wsDe.Range("C1").ColumnWidth = 29.166666666666668
$wsDe.Range("P1").ColumnWidth = 12.833333333333334
